$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3354.5264
$ws.Range("I17").Value = 4763
$ws.Range("J17").Value = 2978.9333
$ws.Range("K17").Value = 14289
$ws.Range("L17").Value = 8936.7999
$ws.Range("M17").Value = -14121
$ws.Range("N17").Value = -9272.7999

$ws.Range("H32").Value = 1436.5555
$ws.Range("J32").Value = 1491.25
$ws.Range("L32").Value = 1491.25
$ws.Range("N32").Value = -2143.25

$ws.Range("H98").Value = 1418.6296
$ws.Range("I98").Value = 1224.619
$ws.Range("K98").Value = 1224.619
$ws.Range("M98").Value = 273.3810000000001

$ws.Range("H106").Value = 2807.2
$ws.Range("I106").Value = 3197.25
$ws.Range("K106").Value = 3197.25
$ws.Range("M106").Value = -2566.25

$ws.Range("H121").Value = 1036.4
$ws.Range("I121").Value = 200
$ws.Range("J121").Value = 1245.5
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 3736.5
$ws.Range("M121").Value = 1147
$ws.Range("N121").Value = -7230.5

$ws.Range("H122").Value = 1418.6296
$ws.Range("I122").Value = 1224.619
$ws.Range("K122").Value = 3673.857
$ws.Range("M122").Value = -1223.857

$ws.Range("H135").Value = 424.76923
$ws.Range("I135").Value = 148.72728
$ws.Range("K135").Value = 1338.54552
$ws.Range("M135").Value = 1196.45448

$ws.Range("H141").Value = 4006759.2
$ws.Range("I141").Value = 5603418.5
$ws.Range("K141").Value = 16810255.5
$ws.Range("M141").Value = -16805075.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4768.213
$ws.Range("I32").Value = 3967.8604
$ws.Range("K32").Value = 3967.8604
$ws.Range("M32").Value = -3680.8604

$ws.Range("H45").Value = 1610.6666
$ws.Range("I45").Value = 1033.75
$ws.Range("K45").Value = 1033.75
$ws.Range("M45").Value = -656.75

$ws.Range("H61").Value = 6284.6523
$ws.Range("I61").Value = 6361.7646
$ws.Range("K61").Value = 6361.7646
$ws.Range("M61").Value = -6149.7646

$ws.Range("H74").Value = 1314.5217
$ws.Range("I74").Value = 457.6111
$ws.Range("K74").Value = 457.6111
$ws.Range("M74").Value = 416.3889

$ws.Range("H77").Value = 1314.5217
$ws.Range("I77").Value = 457.6111
$ws.Range("K77").Value = 2288.0555
$ws.Range("M77").Value = 2079.9445

$ws.Range("H132").Value = 1094.0817
$ws.Range("I132").Value = 987.3555
$ws.Range("K132").Value = 2962.0665
$ws.Range("M132").Value = -432.0664999999999

$ws.Range("H136").Value = 6284.6523
$ws.Range("I136").Value = 6361.7646
$ws.Range("K136").Value = 19085.2938
$ws.Range("M136").Value = -16535.2938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3004.1428
$ws.Range("J20").Value = 3299
$ws.Range("L20").Value = 3299
$ws.Range("N20").Value = -3793

$ws.Range("H22").Value = 868.2
$ws.Range("I22").Value = 844
$ws.Range("J22").Value = 884.3333
$ws.Range("K22").Value = 844
$ws.Range("L22").Value = 884.3333
$ws.Range("M22").Value = -671
$ws.Range("N22").Value = -1230.3333

$ws.Range("H134").Value = 1954.1852
$ws.Range("I134").Value = 1661.35
$ws.Range("K134").Value = 4984.049999999999
$ws.Range("M134").Value = -2449.049999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2124.7273
$ws.Range("I31").Value = 1764.4445
$ws.Range("K31").Value = 1764.4445
$ws.Range("M31").Value = -1469.4445

$ws.Range("H34").Value = 2124.7273
$ws.Range("I34").Value = 1764.4445
$ws.Range("K34").Value = 1764.4445
$ws.Range("M34").Value = -1562.4445

$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

$ws.Range("H58").Value = 6214096
$ws.Range("I58").Value = 14493557
$ws.Range("J58").Value = 4500
$ws.Range("K58").Value = 14493557
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -14493354
$ws.Range("N58").Value = -4906

$ws.Range("H132").Value = 6507
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H134").Value = 1640.1154
$ws.Range("I134").Value = 1617.3334
$ws.Range("J134").Value = 1913.5
$ws.Range("K134").Value = 4852.0002
$ws.Range("L134").Value = 5740.5
$ws.Range("M134").Value = -2317.0002
$ws.Range("N134").Value = -10810.5

$ws.Range("H136").Value = 6214096
$ws.Range("I136").Value = 14493557
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 43480671
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -43478121
$ws.Range("N136").Value = -18600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 289.375
$ws.Range("I2").Value = 278.75
$ws.Range("K2").Value = 1672.5
$ws.Range("M2").Value = -1559.5

$ws.Range("H118").Value = 1164.5
$ws.Range("I118").Value = 609.3333
$ws.Range("J118").Value = 1997.25
$ws.Range("K118").Value = 1827.9999
$ws.Range("L118").Value = 5991.75
$ws.Range("M118").Value = -584.9999
$ws.Range("N118").Value = -8477.75

$ws.Range("H131").Value = 16739.78
$ws.Range("J131").Value = 18133.89
$ws.Range("L131").Value = 54401.67
$ws.Range("N131").Value = -64481.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 8384838
$ws.Range("I21").Value = 25004500
$ws.Range("J21").Value = 75007
$ws.Range("K21").Value = 25004500
$ws.Range("L21").Value = 75007
$ws.Range("M21").Value = -25004327
$ws.Range("N21").Value = -75353

$ws.Range("H30").Value = 8384838
$ws.Range("I30").Value = 25004500
$ws.Range("J30").Value = 75007
$ws.Range("K30").Value = 25004500
$ws.Range("L30").Value = 75007
$ws.Range("M30").Value = -25004395
$ws.Range("N30").Value = -75217

$ws.Range("H102").Value = 1594.1875
$ws.Range("I102").Value = 1275.28
$ws.Range("K102").Value = 1275.28
$ws.Range("M102").Value = 346.72

$ws.Range("H130").Value = 59888
$ws.Range("J130").Value = 59888
$ws.Range("L130").Value = 59888
$ws.Range("N130").Value = -69928

$ws.Range("H132").Value = 2565282.5
$ws.Range("J132").Value = 1850
$ws.Range("L132").Value = 5550
$ws.Range("N132").Value = -10610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 446571.44
$ws.Range("J2").Value = 126000
$ws.Range("L2").Value = 126000
$ws.Range("N2").Value = -126224

$ws.Range("H22").Value = 2120.875
$ws.Range("I22").Value = 1249.5
$ws.Range("J22").Value = 2992.25
$ws.Range("K22").Value = 1249.5
$ws.Range("L22").Value = 2992.25
$ws.Range("M22").Value = -954.5
$ws.Range("N22").Value = -3582.25

$ws.Range("H27").Value = 2120.875
$ws.Range("I27").Value = 1249.5
$ws.Range("J27").Value = 2992.25
$ws.Range("K27").Value = 1249.5
$ws.Range("L27").Value = 2992.25
$ws.Range("M27").Value = -1142.5
$ws.Range("N27").Value = -3206.25

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 10975.529
$ws.Range("I40").Value = 10476.154
$ws.Range("J40").Value = 12598.5
$ws.Range("K40").Value = 10476.154
$ws.Range("L40").Value = 12598.5
$ws.Range("M40").Value = -10340.154
$ws.Range("N40").Value = -12870.5

$ws.Range("H122").Value = 9944.0625
$ws.Range("I122").Value = 10464.091
$ws.Range("K122").Value = 31392.273
$ws.Range("M122").Value = -28942.273

$ws.Range("H136").Value = 2570.9443
$ws.Range("I136").Value = 2447.8333
$ws.Range("K136").Value = 7343.499899999999
$ws.Range("M136").Value = -4793.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 939.8333
$ws.Range("I81").Value = 1037.8
$ws.Range("K81").Value = 2075.6
$ws.Range("M81").Value = -1014.6

$ws.Range("H84").Value = 939.8333
$ws.Range("I84").Value = 1037.8
$ws.Range("K84").Value = 10378
$ws.Range("M84").Value = -5074

$ws.Range("H107").Value = 594.8823
$ws.Range("I107").Value = 422.14285
$ws.Range("K107").Value = 1266.42855
$ws.Range("M107").Value = 653.5714499999999

$ws.Range("H122").Value = 27455.133
$ws.Range("I122").Value = 50199.812
$ws.Range("J122").Value = 1461.2142
$ws.Range("K122").Value = 150599.436
$ws.Range("L122").Value = 4383.642599999999
$ws.Range("M122").Value = -148149.436
$ws.Range("N122").Value = -9283.642599999999

$ws.Range("H132").Value = 1707.4445
$ws.Range("I132").Value = 1280.2593
$ws.Range("J132").Value = 2989
$ws.Range("K132").Value = 3840.7779
$ws.Range("L132").Value = 8967
$ws.Range("M132").Value = -1310.7779
$ws.Range("N132").Value = -14027

$ws.Range("H136").Value = 1081.1025
$ws.Range("I136").Value = 717.0357
$ws.Range("K136").Value = 2151.1071
$ws.Range("M136").Value = 398.8928999999998
